# Change:  <ms>double handful</ms>   ->   <ms>double <bp>handful</bp></ms>
#
# i.e. the run "double handful" is split into "double " + "handful", and
# a new <bp>...</bp> tag pair (in the same Courier-New/blue/9pt style used
# for the other inline markup tags) is wrapped around "handful".
#
# Strategy: locate "double handful" with Find, then clone the character
# formatting of the existing "</ms>" run (which already carries the exact
# Courier-New/blue/sz18/szCs18 rPr we need) via Range.FormattedText onto
# the two insertion points, and finally fix up the literal text of each
# cloned/inserted run. Cloning FormattedText (rather than poking
# individual Font.* properties) is what reproduces the donor run's rPr
# byte-for-byte, including w:cs / w:eastAsia / w:szCs, which plain
# Font.Name/Font.Size assignment does not reliably reproduce here.

$d = $word.ActiveDocument

# --- locate "double handful" ------------------------------------------------
$dh = $d.Content
$found = $dh.Find.Execute("double handful", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'double handful' in the document."
}
$dhStart = $dh.Start
$dhEnd = $dh.End
$handfulStart = $dhStart + 7   # length of "double " == 7

# The run immediately following "double handful" is the existing "</ms>"
# run - our donor for the Courier-New/blue/sz18/szCs18 character style.
$msCloseRng = $d.Range($dhEnd, $dhEnd + 5)
if ($msCloseRng.Text -ne "</ms>") {
    throw "Unexpected text after 'double handful': '$($msCloseRng.Text)'"
}
$donorFormat = $msCloseRng.FormattedText

# --- Step 1: insert "<bp>" right before "handful" ---------------------------
$bpOpenPoint = $d.Range($handfulStart, $handfulStart)
$bpOpenPoint.FormattedText = $donorFormat
$bpOpenInserted = $d.Range($handfulStart, $handfulStart + 5)
$bpOpenInserted.Text = "<bp>"

# --- Step 2: insert "</bp>" right before the (now relocated) "</ms>" --------
$afterHandful = $handfulStart + 4 + 7   # +len("<bp>") +len("handful")
$msCloseRng2 = $d.Range($afterHandful, $afterHandful + 5)
if ($msCloseRng2.Text -ne "</ms>") {
    throw "Unexpected text where '</ms>' was expected: '$($msCloseRng2.Text)'"
}
$bpClosePoint = $d.Range($afterHandful, $afterHandful)
$bpClosePoint.FormattedText = $msCloseRng2.FormattedText
$bpCloseInserted = $d.Range($afterHandful, $afterHandful + 5)
$bpCloseInserted.Text = "</bp>"

Write-Output "Result: $($d.Range($dhStart, $afterHandful + 10).Text)"
